$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = $null
$ws.Range("H48").Value = 1329.1578
$ws.Range("I48").Value = 463.5625
$ws.Range("J48").Value = 5945.6665
$ws.Range("K48").Value = 1390.6875
$ws.Range("L48").Value = 17836.9995
$ws.Range("M48").Value = -1098.6875
$ws.Range("N48").Value = -18420.9995
$ws.Range("H56").Value = 1329.1578
$ws.Range("I56").Value = 463.5625
$ws.Range("J56").Value = 5945.6665
$ws.Range("K56").Value = 1390.6875
$ws.Range("L56").Value = 17836.9995
$ws.Range("M56").Value = -856.6875
$ws.Range("N56").Value = -18904.9995
$ws.Range("H64").Value = 4714.2856
$ws.Range("I64").Value = 3400
$ws.Range("K64").Value = 3400
$ws.Range("M64").Value = -3152
$ws.Range("H67").Value = 4714.2856
$ws.Range("I67").Value = 3400
$ws.Range("K67").Value = 3400
$ws.Range("M67").Value = -2542
$ws.Range("H96").Value = 646.8421
$ws.Range("I96").Value = 373.05884
$ws.Range("K96").Value = 1119.17652
$ws.Range("M96").Value = 253.82348
$ws.Range("H112").Value = 1368.9678
$ws.Range("J112").Value = 1374.4828
$ws.Range("L112").Value = 4123.4484
$ws.Range("N112").Value = -6339.4484
$ws.Range("H132").Value = 2257.6428
$ws.Range("I132").Value = 2017.3334
$ws.Range("K132").Value = 6052.0002
$ws.Range("M132").Value = -3522.0002
$ws.Range("H138").Value = 2431.8438
$ws.Range("J138").Value = 2237.5557
$ws.Range("L138").Value = 6712.6671
$ws.Range("N138").Value = -16992.6671

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 5959.25
$ws.Range("I31").Value = 5959.25
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 5959.25
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -5665.25
$ws.Range("N31").Value = $null
$ws.Range("I33").Value = 10000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 10000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -9671
$ws.Range("N33").Value = $null
$ws.Range("H61").Value = 3211
$ws.Range("I61").Value = 2842.7144
$ws.Range("K61").Value = 2842.7144
$ws.Range("M61").Value = -2630.7144
$ws.Range("H119").Value = 19099.2
$ws.Range("I119").Value = 10000
$ws.Range("J119").Value = 21374
$ws.Range("K119").Value = 10000
$ws.Range("L119").Value = 21374
$ws.Range("M119").Value = -5162
$ws.Range("N119").Value = -31050
$ws.Range("H120").Value = 10000
$ws.Range("J120").Value = 10000
$ws.Range("L120").Value = 10000
$ws.Range("N120").Value = -19676
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("H133").Value = 138261
$ws.Range("J133").Value = 138261
$ws.Range("L133").Value = 138261
$ws.Range("N133").Value = -143321
$ws.Range("H136").Value = 3211
$ws.Range("I136").Value = 2842.7144
$ws.Range("K136").Value = 8528.143199999999
$ws.Range("M136").Value = -5978.143199999999

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = $null
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").Value = $null
$ws.Range("H94").Value = 2303.3684
$ws.Range("I94").Value = 2064.6785
$ws.Range("K94").Value = 2064.6785
$ws.Range("M94").Value = -1613.6785
$ws.Range("H134").Value = 782
$ws.Range("I134").Value = 782
$ws.Range("K134").Value = 2346
$ws.Range("M134").Value = 189

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34763.617
$ws.Range("J31").Value = 27256.5
$ws.Range("L31").Value = 27256.5
$ws.Range("N31").Value = -27846.5
$ws.Range("H34").Value = 34763.617
$ws.Range("J34").Value = 27256.5
$ws.Range("L34").Value = 27256.5
$ws.Range("N34").Value = -27660.5
$ws.Range("H107").Value = 1838.2858
$ws.Range("I107").Value = 1717.0834
$ws.Range("K107").Value = 1717.0834
$ws.Range("M107").Value = 202.9166
$ws.Range("H134").Value = 1470.4286
$ws.Range("I134").Value = 1198.92
$ws.Range("J134").Value = 3733
$ws.Range("K134").Value = 3596.76
$ws.Range("L134").Value = 11199
$ws.Range("M134").Value = -1061.76
$ws.Range("N134").Value = -16269

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 78.588234
$ws.Range("J2").Value = 65.545456
$ws.Range("L2").Value = 393.272736
$ws.Range("N2").Value = -619.272736

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4717.6665
$ws.Range("I70").Value = 4561.143
$ws.Range("J70").Value = 4936.8
$ws.Range("K70").Value = 4561.143
$ws.Range("L70").Value = 4936.8
$ws.Range("M70").Value = -4291.143
$ws.Range("N70").Value = -5476.8
$ws.Range("H73").Value = 4717.6665
$ws.Range("I73").Value = 4561.143
$ws.Range("J73").Value = 4936.8
$ws.Range("K73").Value = 4561.143
$ws.Range("L73").Value = 4936.8
$ws.Range("M73").Value = -3625.143
$ws.Range("N73").Value = -6808.8
$ws.Range("H107").Value = 55563100
$ws.Range("I107").Value = 2399
$ws.Range("J107").Value = 62508188
$ws.Range("K107").Value = 2399
$ws.Range("L107").Value = 62508188
$ws.Range("M107").Value = -479
$ws.Range("N107").Value = -62512028

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4119.769
$ws.Range("I132").Value = 3362.16
$ws.Range("K132").Value = 10086.48
$ws.Range("M132").Value = -7556.48

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1797.8334
$ws.Range("I122").Value = 1757.6
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 5272.799999999999
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = -2822.799999999999
$ws.Range("N122").Value = -10897
$ws.Range("H126").Value = 2217
$ws.Range("I126").Value = 2217
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6651
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4181
$ws.Range("N126").Value = $null
$ws.Range("H132").Value = 2476951.8
$ws.Range("I132").Value = 3961356.5
$ws.Range("K132").Value = 11884069.5
$ws.Range("M132").Value = -11881539.5
